$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.123.97'
$ws.Range("E2").Value = '  +2.97%  '

$ws.Range("D3").Value = '2.656.07'
$ws.Range("E3").Value = '  +2.96%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''594.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.19%  '

$ws.Range("D6").Value = '''156.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.02%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").Value = '''0.594'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.41%  '

$ws.Range("E9").Value = '  +7.79%  '

$ws.Range("E10").Value = '  +4.36%  '

$ws.Range("D11").Value = '''5.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("E12").Value = '  +1.95%  '

$ws.Range("D13").Value = '''29.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.46%  '

$ws.Range("D14").Value = '''0.0000187'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +21.04%  '

$ws.Range("D15").Value = '3.130.71'
$ws.Range("E15").Value = '  +2.79%  '

$ws.Range("D16").Value = '65.012.05'
$ws.Range("E16").Value = '  +3.10%  '

$ws.Range("D17").Value = '2.625.51'
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("E19").Value = '  +1.82%  '

$ws.Range("D20").Value = '''354.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.02%  '

$ws.Range("D21").Value = '''7.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.50%  '

$ws.Range("E22").Value = '  +0.20%  '

$ws.Range("D23").Value = '''68.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.44%  '

$ws.Range("D24").Value = '''1.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.82%  '

$ws.Range("D25").Value = '''9.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.38%  '

$ws.Range("E26").Value = '  -1.04%  '

$ws.Range("D27").Value = '''8.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.62%  '

$ws.Range("E28").Value = '  +2.04%  '

$ws.Range("E29").Value = '  +11.14%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").Value = '''524.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.07%  '

$ws.Range("E32").Value = '  +4.14%  '

$ws.Range("E33").Value = '  +2.24%  '

$ws.Range("D34").Value = '''5.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.03%  '

$ws.Range("D35").Value = '''6.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.62%  '

$ws.Range("E36").Value = '  +3.86%  '

$ws.Range("D37").Value = '''165.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("D38").Value = '''20.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.97%  '

$ws.Range("E39").Value = '  +5.47%  '

$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("D41").Value = '''1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").Value = '''42.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.65%  '

$ws.Range("D43").Value = '''165.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("E44").Value = '  +3.08%  '

$ws.Range("D45").Value = '''0.0617'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.10%  '

$ws.Range("D46").Value = '''22.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("E47").Value = '  +4.89%  '

$ws.Range("D48").Value = '''0.650'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.81%  '

$ws.Range("E49").Value = '  +1.98%  '

$ws.Range("D50").Value = '''0.0987'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.92%  '

$ws.Range("D51").Value = '''19.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.93%  '
